$wb = $excel.ActiveWorkbook

# --- Sheet: "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("B2:B12").NumberFormat = "@"

$ws1.Range("B2").Value = "2024-06-22"
$ws1.Range("C2").Value = "南宁·排球少年ONLY（取消）"
$ws1.Range("D2").Value = "亭洪路45号 水明漾宴会中心"
$ws1.Range("E2").Value = "2024.06.22 09:45-06.22 17:00"
$ws1.Range("F2").Value = 65
$ws1.Range("G2").Value = "不可售"
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=86465"
$ws1.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202405/GaaD97dL1716883956953.jpeg"

$ws1.Range("B3").Value = "2024-07-06"
$ws1.Range("C3").Value = "南宁·小蜜蜂动漫嘉年华2.0"
$ws1.Range("D3").Value = "亭洪路45号 百益上河城"
$ws1.Range("E3").Value = "2024.07.06 10:00-07.06 17:00"
$ws1.Range("F3").Value = 308
$ws1.Range("G3").Value = 50
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=84925"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg"

$ws1.Range("B4").Value = "2024-07-06"
$ws1.Range("C4").Value = "南宁·首届童话梦境Lolita茶会"
$ws1.Range("D4").Value = "明秀东路157号 利泰国际大酒店"
$ws1.Range("E4").Value = "2024.07.06 13:00-07.06 17:00"
$ws1.Range("F4").Value = 164
$ws1.Range("G4").Value = 88
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85776"
$ws1.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg"

$ws1.Range("B5").Value = "2024-07-12"
$ws1.Range("C5").Value = "南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展"
$ws1.Range("D5").Value = "民族大道106号 南宁国际会展中心"
$ws1.Range("E5").Value = "2024.07.12 09:30-07.14 17:00"
$ws1.Range("F5").Value = 190
$ws1.Range("G5").Value = 50
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87182"
$ws1.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202406/x4UZPn301718159475475.jpeg"

$ws1.Range("B6").Value = "2024-07-13"
$ws1.Range("C6").Value = "南宁·0713国乙ONLY"
$ws1.Range("D6").Value = "亭洪路45号 水明漾宴会中心"
$ws1.Range("E6").Value = "2024.07.13 09:30-07.13 21:00"
$ws1.Range("F6").Value = 318
$ws1.Range("G6").Value = 68
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86378"
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg"

$ws1.Range("B7").Value = "2024-07-14"
$ws1.Range("C7").Value = "广西·首届明日方舟only展 - 花庭圣梦"
$ws1.Range("D7").Value = "明秀东路157号 利泰国际大酒店"
$ws1.Range("E7").Value = "2024.07.14 09:00-07.14 18:00"
$ws1.Range("F7").Value = 222
$ws1.Range("G7").Value = "不可售"
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85852"
$ws1.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg"

$ws1.Range("B8").Value = "2024-07-20"
$ws1.Range("C8").Value = "南宁·AB动漫游戏嘉年华"
$ws1.Range("D8").Value = "三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心"
$ws1.Range("E8").Value = "2024.07.20 09:30-07.21 17:00"
$ws1.Range("F8").Value = 2149
$ws1.Range("G8").Value = 60
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84862"
$ws1.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg"

$ws1.Range("B9").Value = "2024-07-20"
$ws1.Range("C9").Value = "横州·第二届海棠动漫游戏嘉年华"
$ws1.Range("D9").Value = "茉莉花大道 横州国际大酒店"
$ws1.Range("E9").Value = "2024.07.20 09:30-07.20 17:00"
$ws1.Range("F9").Value = 371
$ws1.Range("G9").Value = 30
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84799"
$ws1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg"

$ws1.Range("B10").Value = "2024-07-27"
$ws1.Range("C10").Value = "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）"
$ws1.Range("D10").Value = "民族大道106号 南宁国际会展中心"
$ws1.Range("E10").Value = "2024.07.27 09:30-07.28 17:30"
$ws1.Range("F10").Value = 5254
$ws1.Range("G10").Value = 55
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85264"
$ws1.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202405/dZVcS7eE1715155418142.jpeg"

$ws1.Range("B11").Value = "2024-08-03"
$ws1.Range("C11").Value = "南宁·火影忍者only"
$ws1.Range("D11").Value = "厢竹大道65号 桔子酒店"
$ws1.Range("E11").Value = "2024.08.03 10:00-08.03 17:00"
$ws1.Range("F11").Value = 116
$ws1.Range("G11").Value = 68
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=86994"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg"

$ws1.Range("B12").Value = "2024-08-03"
$ws1.Range("C12").Value = "南宁·蔚蓝档案only"
$ws1.Range("D12").Value = "亭洪路45号 百益上河城"
$ws1.Range("E12").Value = "2024.08.03 09:00-08.03 17:00"
$ws1.Range("F12").Value = 353
$ws1.Range("G12").Value = 68
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85370"
$ws1.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg"

$ws1.Rows.Item(13).Delete()

# --- Sheet: "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B2:B15").NumberFormat = "@"

$ws4.Range("B2").Value = "2024-06-22"
$ws4.Range("C2").Value = "南宁·排球少年ONLY（取消）"
$ws4.Range("D2").Value = "亭洪路45号 水明漾宴会中心"
$ws4.Range("E2").Value = "2024.06.22 09:45-06.22 17:00"
$ws4.Range("F2").Value = 65
$ws4.Range("G2").Value = "不可售"
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=86465"
$ws4.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202405/GaaD97dL1716883956953.jpeg"

$ws4.Range("B3").Value = "2024-06-22"
$ws4.Range("C3").Value = "南宁·浪漫古典·百年经典世界名曲音乐会"
$ws4.Range("D3").Value = "广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅"
$ws4.Range("E3").Value = "2024.06.22 20:00-06.22 21:30"
$ws4.Range("F3").Value = 50
$ws4.Range("G3").Value = 135
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=83959"
$ws4.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg"

$ws4.Range("B4").Value = "2024-07-06"
$ws4.Range("C4").Value = "南宁·小蜜蜂动漫嘉年华2.0"
$ws4.Range("D4").Value = "亭洪路45号 百益上河城"
$ws4.Range("E4").Value = "2024.07.06 10:00-07.06 17:00"
$ws4.Range("F4").Value = 308
$ws4.Range("G4").Value = 50
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=84925"
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg"

$ws4.Range("B5").Value = "2024-07-06"
$ws4.Range("C5").Value = "南宁·首届童话梦境Lolita茶会"
$ws4.Range("D5").Value = "明秀东路157号 利泰国际大酒店"
$ws4.Range("E5").Value = "2024.07.06 13:00-07.06 17:00"
$ws4.Range("F5").Value = 164
$ws4.Range("G5").Value = 88
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=85776"
$ws4.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg"

$ws4.Range("B6").Value = "2024-07-12"
$ws4.Range("C6").Value = "南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展"
$ws4.Range("D6").Value = "民族大道106号 南宁国际会展中心"
$ws4.Range("E6").Value = "2024.07.12 09:30-07.14 17:00"
$ws4.Range("F6").Value = 190
$ws4.Range("G6").Value = 50
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87182"
$ws4.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202406/x4UZPn301718159475475.jpeg"

$ws4.Range("B7").Value = "2024-07-13"
$ws4.Range("C7").Value = "南宁·0713国乙ONLY"
$ws4.Range("D7").Value = "亭洪路45号 水明漾宴会中心"
$ws4.Range("E7").Value = "2024.07.13 09:30-07.13 21:00"
$ws4.Range("F7").Value = 318
$ws4.Range("G7").Value = 68
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=86378"
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg"

$ws4.Range("B8").Value = "2024-07-14"
$ws4.Range("C8").Value = "广西·首届明日方舟only展 - 花庭圣梦"
$ws4.Range("D8").Value = "明秀东路157号 利泰国际大酒店"
$ws4.Range("E8").Value = "2024.07.14 09:00-07.14 18:00"
$ws4.Range("F8").Value = 222
$ws4.Range("G8").Value = "不可售"
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=85852"
$ws4.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg"

$ws4.Range("B9").Value = "2024-07-18"
$ws4.Range("C9").Value = "南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》"
$ws4.Range("D9").Value = "龙堤路25号 广西文化艺术中心"
$ws4.Range("E9").Value = "2024.07.18 20:00-07.18 21:30"
$ws4.Range("F9").Value = 9
$ws4.Range("G9").Value = 108
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85816"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg"

$ws4.Range("B10").Value = "2024-07-19"
$ws4.Range("C10").Value = "南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 "
$ws4.Range("D10").Value = "龙堤路25号 广西文化艺术中心"
$ws4.Range("E10").Value = "2024.07.19 20:00-07.19 22:00"
$ws4.Range("F10").Value = 14
$ws4.Range("G10").Value = 108
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85831"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg"

$ws4.Range("B11").Value = "2024-07-20"
$ws4.Range("C11").Value = "南宁·AB动漫游戏嘉年华"
$ws4.Range("D11").Value = "三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心"
$ws4.Range("E11").Value = "2024.07.20 09:30-07.21 17:00"
$ws4.Range("F11").Value = 2149
$ws4.Range("G11").Value = 60
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=84862"
$ws4.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg"

$ws4.Range("B12").Value = "2024-07-20"
$ws4.Range("C12").Value = "横州·第二届海棠动漫游戏嘉年华"
$ws4.Range("D12").Value = "茉莉花大道 横州国际大酒店"
$ws4.Range("E12").Value = "2024.07.20 09:30-07.20 17:00"
$ws4.Range("F12").Value = 371
$ws4.Range("G12").Value = 30
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=84799"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg"

$ws4.Range("B13").Value = "2024-07-27"
$ws4.Range("C13").Value = "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）"
$ws4.Range("D13").Value = "民族大道106号 南宁国际会展中心"
$ws4.Range("E13").Value = "2024.07.27 09:30-07.28 17:30"
$ws4.Range("F13").Value = 5254
$ws4.Range("G13").Value = 55
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85264"
$ws4.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202405/dZVcS7eE1715155418142.jpeg"

$ws4.Range("B14").Value = "2024-08-03"
$ws4.Range("C14").Value = "南宁·火影忍者only"
$ws4.Range("D14").Value = "厢竹大道65号 桔子酒店"
$ws4.Range("E14").Value = "2024.08.03 10:00-08.03 17:00"
$ws4.Range("F14").Value = 116
$ws4.Range("G14").Value = 68
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86994"
$ws4.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg"

$ws4.Range("B15").Value = "2024-08-03"
$ws4.Range("C15").Value = "南宁·蔚蓝档案only"
$ws4.Range("D15").Value = "亭洪路45号 百益上河城"
$ws4.Range("E15").Value = "2024.08.03 09:00-08.03 17:00"
$ws4.Range("F15").Value = 353
$ws4.Range("G15").Value = 68
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=85370"
$ws4.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg"

$ws4.Rows.Item(16).Delete()
